$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.787.93"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.644.36"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("D5").Value = "'216.82"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("D9").Value = "'0.0626"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("D11").Value = "'0.0842"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "1.867.13"
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").Value = "1.646.71"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").Value = "'4.18"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "'0.527"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("D16").Value = "'64.68"
$ws.Range("E16").Value = "  -3.30%  "
$ws.Range("D17").Value = "26.773.73"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("E18").Value = "  -2.49%  "
$ws.Range("D19").Value = "'214.36"
$ws.Range("E19").Value = "  -3.15%  "
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D21").Value = "'4.36"
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("E22").Value = "  +12.43%  "
$ws.Range("E23").Value = "  -1.24%  "
$ws.Range("D24").Value = "'9.37"
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("D25").Value = "'144.97"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("D27").Value = "'0.118"
$ws.Range("E27").Value = "  -2.30%  "
$ws.Range("D28").Value = "'7.10"
$ws.Range("D29").Value = "'15.71"
$ws.Range("E29").Value = "  -1.59%  "
$ws.Range("D30").Value = "'0.0515"
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").Value = "'3.33"
$ws.Range("E32").Value = "  -3.29%  "
$ws.Range("E33").Value = "  -2.04%  "
$ws.Range("D34").Value = "1.288.10"
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("E35").Value = "  -1.81%  "
$ws.Range("D36").Value = "'2.44"
$ws.Range("E36").Value = "  +1.24%  "
$ws.Range("E37").Value = "  -4.15%  "
$ws.Range("E38").Value = "  +2.08%  "
$ws.Range("D39").Value = "'0.827"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("E41").Value = "  -1.01%  "
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").Value = "'5.36"
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("D44").Value = "1.793.28"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "'91.58"
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("D46").Value = "'60.06"
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("D49").Value = "'0.0522"
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("D50").Value = "'7.72"
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("E51").Value = "  -0.33%  "
